$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated Price (D) and Volume(1h) (E) values for the refreshed crypto symbol list.
# Values are kept as literal text (NumberFormat "@") to match the source data, which stores
# these figures as plain strings rather than numeric/percentage values.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '261.17'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '1.77%'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '1.46%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '4.707'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-1.20%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.06089'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '2.83%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.677'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '1.00%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8460'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-0.60%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9271'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '0.21%'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '2.05%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.04777'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '14.58%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07102'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '1.20%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03084'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '1.03%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09067'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001534'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '0.41%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0006077'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-94.09%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.006122'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '0.96%'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-0.66%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.149'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-0.50%'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.183'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-0.72%'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '2.33%'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '0.20%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.084'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '4.80%'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '0.00%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001222'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '0.35%'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '4.89%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001200'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '0.12%'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '3.50%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03875'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '2.48%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.004080'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-34.87%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.01626'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '15.30%'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '0.86%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005138'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-3.75%'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '0.11%'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-43.90%'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '23.86%'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '0.11%'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.11%'
